$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(6, 20, 4, 0),
    @(4, 15, 5, 5),
    @(2, 8, 7, 12),
    @(4, 5, 3, 15),
    @(4, 5, 3, 15),
    @(3, 16, 5, 4),
    @(6, 2, 3, 18),
    @(6, 13, 4, 7),
    @(2, 5, 5, 15),
    @(3, 15, 4, 5),
    @(2, 2, 3, 18),
    @(5, 17, 4, 3),
    @(4, 8, 5, 12),
    @(3, 12, 4, 8)
)

$startRow = 1210
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}

$ws.Range("A1224").Select()
$excel.ActiveWindow.ScrollRow = 1217
$excel.ActiveWindow.ScrollColumn = 1
